$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.137.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.25%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.554.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.59%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'549.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -4.06%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'140.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.10%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.21%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.589"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.00%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.558.41"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.08%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.87%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -2.24%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +5.08%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +1.81%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.000.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.72%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'59.122.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.22%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'22.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.75%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -2.11%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.558.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.09%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.67%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'335.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.00%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.04%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.34%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.993"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.46%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +5.11%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'62.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -5.20%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.18%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -3.41%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.75%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0757"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.43%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.06%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.06%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'6.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.91%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'158.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.20%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'18.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.39%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'4.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.79%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.06%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.886"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.94%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'37.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.89%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.838"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -5.44%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'3.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.84%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'283.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.33%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'134.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +6.50%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.30%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.92%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'10.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.41%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.587"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.52%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -3.23%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0232"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'1.944.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.27%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'18.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.90%  "
$ws.Range("E51").Style = "Normal"
